# Correct the misspelled "EtoH" treatment label to the proper chemical
# notation "EtOH" (ethanol) throughout the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("EtoH", "EtOH")
